$wb = $excel.ActiveWorkbook

# --- Parameter sheets: est (B) / se (C) for rows 2 and 3 ---

$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.97605246674702
$ws.Range("C2").Value = 0.323540867216825
$ws.Range("B3").Value = 0.137278719565766
$ws.Range("C3").Value = 0.216077907045135

$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.4381676875878
$ws.Range("C2").Value = 0.484812865854487
$ws.Range("B3").Value = -1.07201620635763
$ws.Range("C3").Value = 0.192364462982496

$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.28588536418224
$ws.Range("C2").Value = 0.118148159030896
$ws.Range("B3").Value = 1.97834259892908
$ws.Range("C3").Value = 0.38226490490551

$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.65027185824501
$ws.Range("C2").Value = 0.225090686544683
$ws.Range("B3").Value = 0.00000447698641238702
$ws.Range("C3").Value = 0.0270325818460841

# --- Covariance-matrix sheets ---

$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.104678692759415
$ws.Range("B2").Value = -0.0591706299533721
$ws.Range("A3").Value = -0.0591706299533721
$ws.Range("B3").Value = 0.0466896619130061

$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.235043514898041
$ws.Range("B2").Value = -0.0878022120874772
$ws.Range("A3").Value = -0.0878022120874772
$ws.Range("B3").Value = 0.0370040866185441

$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.01395898748239
$ws.Range("B2").Value = -0.00159880999586965
$ws.Range("A3").Value = -0.00159880999586965
$ws.Range("B3").Value = 0.146126457522419

$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0506658171691567
$ws.Range("B2").Value = -0.00415997704233024
$ws.Range("A3").Value = -0.00415997704233024
$ws.Range("B3").Value = 0.000730760481265237

$wb.Save()
